$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix inconsistent "Feature" column text for the LookupValue UT / Posting section
# (rows were split between "LookupValue" and "Lookup Value" instead of the
# canonical "LookupValue UT" used by every other row in the table).
$ws.Range("A11:A18").Value = "LookupValue UT"

# Leave the cursor where the author last left it after fixing the values.
$ws.Range("A15").Select()
